$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D73").Value = 45007
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 60
$ws.Range("N73").Value = 7000
$ws.Range("O73").Value = 7500
$ws.Range("P73").Value = 7250
$ws.Range("R73").Value = "Perú"
$ws.Range("S73").Value = 1812

$ws.Range("D74").Value = 44959
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 60
$ws.Range("N74").Value = 6500
$ws.Range("O74").Value = 7000
$ws.Range("P74").Value = 6750
$ws.Range("R74").Value = "Perú"
$ws.Range("S74").Value = 1688

$ws.Range("D75").Value = 44841
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 80
$ws.Range("N75").Value = 7500
$ws.Range("O75").Value = 8000
$ws.Range("P75").Value = 7750
$ws.Range("R75").Value = "Brasil"
$ws.Range("S75").Value = 1938

$ws.Range("D76").Value = 44196
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 50
$ws.Range("N76").Value = 6000
$ws.Range("O76").Value = 6000
$ws.Range("P76").Value = 6000
$ws.Range("R76").Value = "Perú"
$ws.Range("S76").Value = 1500

$ws.Range("D77").Value = 44232
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 60
$ws.Range("N77").Value = 6000
$ws.Range("O77").Value = 6000
$ws.Range("P77").Value = 6000
$ws.Range("R77").Value = "Perú"
$ws.Range("S77").Value = 1500

$ws.Range("D78").Value = 44475
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 60
$ws.Range("N78").Value = 8500
$ws.Range("O78").Value = 9000
$ws.Range("P78").Value = 8750
$ws.Range("R78").Value = "Brasil"
$ws.Range("S78").Value = 2188

$ws.Range("D79").Value = 44181
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 50
$ws.Range("N79").Value = 6000
$ws.Range("O79").Value = 6000
$ws.Range("P79").Value = 6000
$ws.Range("R79").Value = "Perú"
$ws.Range("S79").Value = 1500

$ws.Range("D80").Value = 44181
$ws.Range("L80").Value = "Segunda"
$ws.Range("M80").Value = 40
$ws.Range("N80").Value = 5000
$ws.Range("O80").Value = 5000
$ws.Range("P80").Value = 5000
$ws.Range("R80").Value = "Perú"
$ws.Range("S80").Value = 1250

$ws.Range("D81").Value = 44943
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 60
$ws.Range("N81").Value = 7000
$ws.Range("O81").Value = 7000
$ws.Range("P81").Value = 7000
$ws.Range("R81").Value = "Perú"
$ws.Range("S81").Value = 1750

$ws.Range("D82").Value = 44420
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 40
$ws.Range("N82").Value = 8500
$ws.Range("O82").Value = 9000
$ws.Range("P82").Value = 8750
$ws.Range("R82").Value = "Perú"
$ws.Range("S82").Value = 2188

$ws.Range("D83").Value = 44419
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 60
$ws.Range("N83").Value = 8500
$ws.Range("O83").Value = 9000
$ws.Range("P83").Value = 8750
$ws.Range("R83").Value = "Perú"
$ws.Range("S83").Value = 2188

$ws.Range("D84").Value = 44186
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 40
$ws.Range("N84").Value = 6000
$ws.Range("O84").Value = 6000
$ws.Range("P84").Value = 6000
$ws.Range("R84").Value = "Perú"
$ws.Range("S84").Value = 1500

$ws.Range("D85").Value = 44222
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 50
$ws.Range("N85").Value = 6000
$ws.Range("O85").Value = 6000
$ws.Range("P85").Value = 6000
$ws.Range("R85").Value = "Perú"
$ws.Range("S85").Value = 1500

$ws.Range("D86").Value = 44433
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 60
$ws.Range("N86").Value = 8500
$ws.Range("O86").Value = 9000
$ws.Range("P86").Value = 8750
$ws.Range("R86").Value = "Perú"
$ws.Range("S86").Value = 2188

$ws.Range("D87").Value = 44455
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 60
$ws.Range("N87").Value = 8500
$ws.Range("O87").Value = 9000
$ws.Range("P87").Value = 8750
$ws.Range("R87").Value = "Brasil"
$ws.Range("S87").Value = 2188

$ws.Range("D88").Value = 44764
$ws.Range("L88").Value = "Primera"
$ws.Range("M88").Value = 30
$ws.Range("N88").Value = 8500
$ws.Range("O88").Value = 9000
$ws.Range("P88").Value = 8750
$ws.Range("R88").Value = "Brasil"
$ws.Range("S88").Value = 2188

$ws.Range("D89").Value = 44974
$ws.Range("L89").Value = "Primera"
$ws.Range("M89").Value = 50
$ws.Range("N89").Value = 7500
$ws.Range("O89").Value = 7500
$ws.Range("P89").Value = 7500
$ws.Range("R89").Value = "Perú"
$ws.Range("S89").Value = 1875

$ws.Range("D90").Value = 44195
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 50
$ws.Range("N90").Value = 6000
$ws.Range("O90").Value = 6000
$ws.Range("P90").Value = 6000
$ws.Range("R90").Value = "Perú"
$ws.Range("S90").Value = 1500

$ws.Range("D91").Value = 44435
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 60
$ws.Range("N91").Value = 8500
$ws.Range("O91").Value = 9000
$ws.Range("P91").Value = 8750
$ws.Range("R91").Value = "Perú"
$ws.Range("S91").Value = 2188

$ws.Range("D92").Value = 45005
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 50
$ws.Range("N92").Value = 7000
$ws.Range("O92").Value = 7000
$ws.Range("P92").Value = 7000
$ws.Range("R92").Value = "Perú"
$ws.Range("S92").Value = 1750

$ws.Range("D93").Value = 44603
$ws.Range("L93").Value = "Primera"
$ws.Range("M93").Value = 60
$ws.Range("N93").Value = 7000
$ws.Range("O93").Value = 7500
$ws.Range("P93").Value = 7250
$ws.Range("R93").Value = "Perú"
$ws.Range("S93").Value = 1812

$ws.Range("D94").Value = 44650
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 120
$ws.Range("N94").Value = 7500
$ws.Range("O94").Value = 8000
$ws.Range("P94").Value = 7750
$ws.Range("R94").Value = "Perú"
$ws.Range("S94").Value = 1938

$ws.Range("D95").Value = 44229
$ws.Range("L95").Value = "Especial"
$ws.Range("M95").Value = 25
$ws.Range("N95").Value = 6000
$ws.Range("O95").Value = 6000
$ws.Range("P95").Value = 6000
$ws.Range("R95").Value = "Perú"
$ws.Range("S95").Value = 1500

$ws.Range("D96").Value = 44229
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 35
$ws.Range("N96").Value = 6000
$ws.Range("O96").Value = 6000
$ws.Range("P96").Value = 6000
$ws.Range("R96").Value = "Perú"
$ws.Range("S96").Value = 1500

$ws.Range("D97").Value = 44292
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 30
$ws.Range("N97").Value = 7500
$ws.Range("O97").Value = 8000
$ws.Range("P97").Value = 7750
$ws.Range("R97").Value = "Perú"
$ws.Range("S97").Value = 1938

$ws.Range("D98").Value = 44868
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 120
$ws.Range("N98").Value = 8000
$ws.Range("O98").Value = 8500
$ws.Range("P98").Value = 8250
$ws.Range("R98").Value = "Brasil"
$ws.Range("S98").Value = 2062

$ws.Range("D99").Value = 44253
$ws.Range("L99").Value = "Especial"
$ws.Range("M99").Value = 30
$ws.Range("N99").Value = 6000
$ws.Range("O99").Value = 6000
$ws.Range("P99").Value = 6000
$ws.Range("R99").Value = "Perú"
$ws.Range("S99").Value = 1500

$ws.Range("D100").Value = 44253
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 50
$ws.Range("N100").Value = 6000
$ws.Range("O100").Value = 6000
$ws.Range("P100").Value = 6000
$ws.Range("R100").Value = "Perú"
$ws.Range("S100").Value = 1500

$ws.Range("D101").Value = 44998
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 40
$ws.Range("N101").Value = 8000
$ws.Range("O101").Value = 8000
$ws.Range("P101").Value = 8000
$ws.Range("R101").Value = "Perú"
$ws.Range("S101").Value = 2000

$ws.Range("D102").Value = 44473
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 60
$ws.Range("N102").Value = 8500
$ws.Range("O102").Value = 9000
$ws.Range("P102").Value = 8750
$ws.Range("R102").Value = "Brasil"
$ws.Range("S102").Value = 2188

$ws.Range("D103").Value = 44874
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 60
$ws.Range("N103").Value = 9000
$ws.Range("O103").Value = 10000
$ws.Range("P103").Value = 9500
$ws.Range("R103").Value = "Brasil"
$ws.Range("S103").Value = 2375

$ws.Range("D104").Value = 44426
$ws.Range("L104").Value = "Primera"
$ws.Range("M104").Value = 40
$ws.Range("N104").Value = 8500
$ws.Range("O104").Value = 9000
$ws.Range("P104").Value = 8750
$ws.Range("R104").Value = "Perú"
$ws.Range("S104").Value = 2188

$ws.Range("D105").Value = 45008
$ws.Range("L105").Value = "Especial"
$ws.Range("M105").Value = 150
$ws.Range("N105").Value = 8000
$ws.Range("O105").Value = 8000
$ws.Range("P105").Value = 8000
$ws.Range("R105").Value = "Perú"
$ws.Range("S105").Value = 2000

$ws.Range("D106").Value = 45008
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 60
$ws.Range("N106").Value = 7000
$ws.Range("O106").Value = 7000
$ws.Range("P106").Value = 7000
$ws.Range("R106").Value = "Perú"
$ws.Range("S106").Value = 1750

$ws.Range("D107").Value = 44209
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 50
$ws.Range("N107").Value = 6000
$ws.Range("O107").Value = 6000
$ws.Range("P107").Value = 6000
$ws.Range("R107").Value = "Perú"
$ws.Range("S107").Value = 1500

$ws.Range("A108").Value = 7
$ws.Range("B108").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C108").Value = "Ñuble"
$ws.Range("D108").Value = 44210
$ws.Range("D108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E108").Value = 16
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100108
$ws.Range("H108").Value = "Tropicales y subtropicales"
$ws.Range("I108").Value = 100108002
$ws.Range("J108").Value = "Mango"
$ws.Range("K108").Value = "Sin especificar"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 40
$ws.Range("N108").Value = 6000
$ws.Range("O108").Value = 6000
$ws.Range("P108").Value = 6000
$ws.Range("Q108").Value = "$/bandeja 4 kilos"
$ws.Range("R108").Value = "Perú"
$ws.Range("S108").Value = 1500
$ws.Range("T108").Value = 4

$ws.Range("A109").Value = 7
$ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C109").Value = "Ñuble"
$ws.Range("D109").Value = 44189
$ws.Range("D109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E109").Value = 16
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100108
$ws.Range("H109").Value = "Tropicales y subtropicales"
$ws.Range("I109").Value = 100108002
$ws.Range("J109").Value = "Mango"
$ws.Range("K109").Value = "Sin especificar"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 40
$ws.Range("N109").Value = 6000
$ws.Range("O109").Value = 6000
$ws.Range("P109").Value = 6000
$ws.Range("Q109").Value = "$/bandeja 4 kilos"
$ws.Range("R109").Value = "Perú"
$ws.Range("S109").Value = 1500
$ws.Range("T109").Value = 4

$ws.Range("A110").Value = 7
$ws.Range("B110").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C110").Value = "Ñuble"
$ws.Range("D110").Value = 44601
$ws.Range("D110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E110").Value = 16
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100108
$ws.Range("H110").Value = "Tropicales y subtropicales"
$ws.Range("I110").Value = 100108002
$ws.Range("J110").Value = "Mango"
$ws.Range("K110").Value = "Sin especificar"
$ws.Range("L110").Value = "Primera"
$ws.Range("M110").Value = 60
$ws.Range("N110").Value = 7000
$ws.Range("O110").Value = 7500
$ws.Range("P110").Value = 7250
$ws.Range("Q110").Value = "$/bandeja 4 kilos"
$ws.Range("R110").Value = "Perú"
$ws.Range("S110").Value = 1812
$ws.Range("T110").Value = 4
